## Ethan Drost Resume — targeted corrections
## 1. Phone number typo: (506) 459-8314 -> (506) 259-8314
## 2. Company name typo: Griffen endeavors ltd. -> Griffin endeavors ltd.
## 3. Degree typo: "electracal" -> "electrical" (spelled electr + I + cal, per source edit)
##    and tidy up the stray proofing spacing/markers around "engineering,"
## 4. Skills punctuation: "Java. Python" -> "Java, Python"

$d = $word.ActiveDocument

$wdFindContinue = 1
$wdReplaceOne   = 1

# --- 1) Phone number -------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("459-8314", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "259-8314", $wdReplaceOne) | Out-Null

# --- 2) Griffen -> Griffin --------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Griffen endeavors ltd.", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "Griffin endeavors ltd.", $wdReplaceOne) | Out-Null

# --- 3) Bacholer of science in electracal ... engineering, -----------
$rng = $d.Content
$rng.Find.Execute("electracal", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "electrIcal", $wdReplaceOne) | Out-Null

# Re-touch the "computer engineering,<spaces>" span so stray proofing
# marks that Word's grammar checker had attached to the old wording are
# cleared out (consistent with the whole phrase being re-validated after
# the edit above), while leaving the single trailing space that precedes
# "The university..." alone.
$rng = $d.Content
$rng.Find.Execute("computer engineering,          ", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "computer engineering,          ", $wdReplaceOne) | Out-Null

# --- 4) Programming:  Java. Python, SQL, MATLAB, C --------------------
$rng = $d.Content
$rng.Find.Execute("Java. Python, SQL, MATLAB, C", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "Java, Python, SQL, MATLAB, C", $wdReplaceOne) | Out-Null
